# Updated test data for normal load, cable capacitance etc
#
# 1) Duplicate the "Add Panels" sheet (as it stood before the value edits
#    below) to a new trailing sheet named "Sheet2" - this is what shows up
#    in the diff as the brand-new xl/worksheets/sheet3.xml part.
# 2) On the original "Add Panels" sheet, update the Battery Standby test
#    data (F8: Stand By Hours-driven figure, K8: Minimum Battery size) and
#    leave the selection on K9, matching the recorded cursor move.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Add Panels")

# Copy "Add Panels" to the end of the workbook, then rename the copy.
$src.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$copy = $wb.Worksheets.Item($wb.Worksheets.Count)
$copy.Name = "Sheet2"

# The duplicated sheet ends up fully selected (A1:XFD1048576) and not the
# active tab once we flip back to "Add Panels" below.
$copy.Range("A1:XFD1048576").Select()

# Go back to "Add Panels" and apply the updated test values.
$src.Activate()
$src.Range("F8").Value = 8.6
$src.Range("K8").Value = 37.86
$src.Range("K9").Select()
